# "more bulk upload fixes"
#
# Adds a new "Date Created (Year)*" column (E) to Sheet1 of the bulk-upload
# template: a header cell in E1, a sample value (2000) in E2, and two blank
# formatted cells below it (E3:E4) matching the header/value font styling
# (explicit black font color) so the column reads as fully formatted through
# row 4.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New header + sample value in column E.
$ws.Range("E1").Value = "Date Created (Year)*"
$ws.Range("E2").Value = 2000

# Give the whole new column (header through the extra blank rows) the same
# explicit black font color, which mints the new font/cell style used by
# E1:E4.
$ws.Range("E1:E4").Font.Color = 0

# Match the author's selection left on the sheet after the edit.
$ws.Range("E3:E4").Select()
